$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$newValue = "<p>中国工商银行股份有限公司兰州中山支行2703050119200044428</p><p>兰州优行驰辰网络科技有限公司</p><p>91620100MA72R6U81E甘肃省兰州市安宁区北滨河西路530号连铝大厦19楼0931-8617788</p><p>招商银行兰州中央广场支行931905025010909</p><p>031/6&gt;631-+/61**//7&lt;/*101*02+1&gt;9275+&lt;6-615*981083&gt;&lt;7+3-&gt;9491558-73&lt;1-1816-98340298834-4/745*4-01&gt;+78198-+/96301/</p><p>潘佳昕王玉平吴艳</p><p>*运输服务*客运服务费次1279.2279.20免税***</p><p>*运输服务*客运服务费-6.67免税***</p><p>¥272.53***</p><p>贰佰柒拾贰圆伍角叁分¥272.53</p>"

$ws.Range("A2").Value = $newValue
